$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3811094136979136
$ws.Range("C2").Value = 0.1923305032046585
$ws.Range("E2").Value = 0.4029314417661851
$ws.Range("F2").Value = 0.4443680307746263
$ws.Range("G2").Value = 0.002501565051071691
$ws.Range("I2").Value = 1.856491428754794
$ws.Range("K2").Value = 0.5716404330439957
$ws.Range("B3").Value = 0.3536534848279587
$ws.Range("C3").Value = 0.1742865515051051
$ws.Range("E3").Value = 0.3514740142551886
$ws.Range("F3").Value = 0.387822817061874
$ws.Range("G3").Value = 0.002506323312352612
$ws.Range("I3").Value = 1.77356445341519
$ws.Range("K3").Value = 0.5260595393841072
$ws.Range("B4").Value = 0.3371396458981337
$ws.Range("C4").Value = 0.163294212608605
$ws.Range("E4").Value = 0.3200003049035729
$ws.Range("F4").Value = 0.3531389305168915
$ws.Range("G4").Value = 0.002509390414907902
$ws.Range("I4").Value = 1.722877816663868
$ws.Range("K4").Value = 0.4984940867990701
$ws.Range("B5").Value = 0.3304960255400715
$ws.Range("C5").Value = 0.1588359590628556
$ws.Range("E5").Value = 0.3072025220461256
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002510677014384208
$ws.Range("I5").Value = 1.702277617992067
$ws.Range("K5").Value = 0.4873656761031384
$ws.Range("B6").Value = 0.3293980270882173
$ws.Range("C6").Value = 0.1580969335401221
$ws.Range("E6").Value = 0.3050790789481113
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002510892876051363
$ws.Range("I6").Value = 1.698860205058196
$ws.Range("K6").Value = 0.4855241022968073
$ws.Range("B7").Value = 0.3370497007462632
$ws.Range("C7").Value = 0.1632340019241383
$ws.Range("E7").Value = 0.3198275993852064
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002509407617410132
$ws.Range("I7").Value = 1.722599776437292
$ws.Range("K7").Value = 0.498343582702887
$ws.Range("B8").Value = 0.3715707609107142
$ws.Range("C8").Value = 0.1860905856852355
$ws.Range("E8").Value = 0.3851622005618793
$ws.Range("F8").Value = 0.4248636149813336
$ws.Range("G8").Value = 0.002503175583972749
$ws.Range("I8").Value = 1.82784844107303
$ws.Range("K8").Value = 0.5558359170362053
$ws.Range("B9").Value = 0.4420318647516126
$ws.Range("C9").Value = 0.2316295693869961
$ws.Range("E9").Value = 0.5143659878155944
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.00249210268794428
$ws.Range("I9").Value = 2.03621874455861
$ws.Range("K9").Value = 0.6719845832574549
$ws.Range("B10").Value = 0.495539919037185
$ws.Range("C10").Value = 0.2655683545081331
$ws.Range("E10").Value = 0.6101280498869528
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002484658252580767
$ws.Range("I10").Value = 2.190730168883363
$ws.Range("K10").Value = 0.7594930758648104
$ws.Range("B11").Value = 0.5202724950116533
$ws.Range("C11").Value = 0.2811223345443352
$ws.Range("E11").Value = 0.6539133705943954
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002481419638734749
$ws.Range("I11").Value = 2.261377649292513
$ws.Range("K11").Value = 0.7997971496401703
$ws.Range("B12").Value = 0.5296952960759711
$ws.Range("C12").Value = 0.2870295103787157
$ws.Range("E12").Value = 0.6705288200136579
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002480214377281954
$ws.Range("I12").Value = 2.288185279109854
$ws.Range("K12").Value = 0.8151322284107039
$ws.Range("B13").Value = 0.5276633739867975
$ws.Range("C13").Value = 0.2857565175613388
$ws.Range("E13").Value = 0.6669487871793933
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002480473014270105
$ws.Range("I13").Value = 2.282409279533255
$ws.Range("K13").Value = 0.8118262831338541
$ws.Range("B14").Value = 0.5210465651741458
$ws.Range("C14").Value = 0.2816079724015026
$ws.Range("E14").Value = 0.6552796175186018
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002481320058354454
$ws.Range("I14").Value = 2.263582008770982
$ws.Range("K14").Value = 0.8010573081951122
$ws.Range("B15").Value = 0.5170010409094061
$ws.Range("C15").Value = 0.2790691316775167
$ws.Range("E15").Value = 0.6481365453332728
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002481841645262937
$ws.Range("I15").Value = 2.252057017066448
$ws.Range("K15").Value = 0.7944705252550364
$ws.Range("B16").Value = 0.4939315339903203
$ws.Range("C16").Value = 0.2645542380988104
$ws.Range("E16").Value = 0.607271301868991
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002484872867875275
$ws.Range("I16").Value = 2.186120702290339
$ws.Range("K16").Value = 0.7568692162090542
$ws.Range("B17").Value = 0.4798799817282031
$ws.Range("C17").Value = 0.2556797671778384
$ws.Range("E17").Value = 0.5822608629553088
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002486770208109688
$ws.Range("I17").Value = 2.145765295432625
$ws.Range("K17").Value = 0.7339301110158658
$ws.Range("B18").Value = 0.4718346995150284
$ws.Range("C18").Value = 0.2505861947437324
$ws.Range("E18").Value = 0.567896304251633
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002487875436760244
$ws.Range("I18").Value = 2.122587569573199
$ws.Range("K18").Value = 0.7207827420138528
$ws.Range("B19").Value = 0.4691170010247561
$ws.Range("C19").Value = 0.2488634311647502
$ws.Range("E19").Value = 0.5630361943961617
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002488252044643646
$ws.Range("I19").Value = 2.114745667690798
$ws.Range("K19").Value = 0.7163392276464151
$ws.Range("B20").Value = 0.4813719790739412
$ws.Range("C20").Value = 0.2566233483264568
$ws.Range("E20").Value = 0.5849210951020893
$ws.Range("F20").Value = 0.6429339538360921
$ws.Range("G20").Value = 0.002486566792353711
$ws.Range("I20").Value = 2.15005769247756
$ws.Range("K20").Value = 0.7363671846928526
$ws.Range("B21").Value = 0.5229885257128046
$ws.Range("C21").Value = 0.2828260278526784
$ws.Range("E21").Value = 0.6587061651737969
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002481070688687407
$ws.Range("I21").Value = 2.269110516706775
$ws.Range("K21").Value = 0.8042184294636741
$ws.Range("B22").Value = 0.5505206291365425
$ws.Range("C22").Value = 0.3000516940417128
$ws.Range("E22").Value = 0.7071337642704805
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002477601773314077
$ws.Range("I22").Value = 2.347240318724914
$ws.Range("K22").Value = 0.8489881886193302
$ws.Range("B23").Value = 0.5357954560720088
$ws.Range("C23").Value = 0.2908485912146261
$ws.Range("E23").Value = 0.6812673400948484
$ws.Range("F23").Value = 0.7472568307916134
$ws.Range("G23").Value = 0.002479441979502559
$ws.Range("I23").Value = 2.305510431209768
$ws.Range("K23").Value = 0.8250543372894583
$ws.Range("B24").Value = 0.4806973439043531
$ws.Range("C24").Value = 0.2561967289297513
$ws.Range("E24").Value = 0.5837183598996347
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002486658711755444
$ws.Range("I24").Value = 2.148117027959927
$ws.Range("K24").Value = 0.7352652572466809
$ws.Range("B25").Value = 0.4226686731077791
$ws.Range("C25").Value = 0.2192287656314136
$ws.Range("E25").Value = 0.4792790027021141
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.002494976228167505
$ws.Range("I25").Value = 1.979615053197378
$ws.Range("K25").Value = 0.640188665000494
